$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A73").Value2 = 71
$ws.Range("B73").Value2 = 7646750
$ws.Range("C73").Value2 = 'Australia ALeague'
$ws.Range("D73").Value2 = 'Australia ALeague'
$ws.Range("E73").Value2 = 45305.23958333334
$ws.Range("F73").Value2 = 'Perth Glory'
$ws.Range("G73").Value2 = 'Wellington Phoenix'
$ws.Range("H73").Value2 = 3
$ws.Range("I73").Value2 = 4
$ws.Range("J73").Value2 = 'A'
$ws.Range("K73").Value2 = 2.45
$ws.Range("L73").Value2 = 3.75
$ws.Range("M73").Value2 = 2.55
$ws.Range("N73").Value2 = 3.1
$ws.Range("O73").Value2 = 3.8
$ws.Range("P73").Value2 = 2.05
$ws.Range("Q73").Value2 = 0.25
$ws.Range("R73").Value2 = 2
$ws.Range("S73").Value2 = 1.85
$ws.Range("T73").Value2 = 3
$ws.Range("U73").Value2 = 1.925
$ws.Range("V73").Value2 = 1.925
$ws.Range("W73").Value2 = -1
$ws.Range("X73").Value2 = -1
$ws.Range("Y73").Value2 = 1.05
$ws.Range("Z73").Value2 = -1
$ws.Range("AA73").Value2 = 0.8500000000000001
$ws.Range("AB73").Value2 = 0.925
$ws.Range("AC73").Value2 = -1
$ws.Range("A74").Value2 = 72
$ws.Range("B74").Value2 = 7646749
$ws.Range("C74").Value2 = 'Australia ALeague'
$ws.Range("D74").Value2 = 'Australia ALeague'
$ws.Range("E74").Value2 = 45305.23958333334
$ws.Range("F74").Value2 = 'Brisbane Roar'
$ws.Range("G74").Value2 = 'Newcastle Jets'
$ws.Range("H74").Value2 = 3
$ws.Range("I74").Value2 = 2
$ws.Range("J74").Value2 = 'H'
$ws.Range("K74").Value2 = 1.909
$ws.Range("L74").Value2 = 4
$ws.Range("M74").Value2 = 3.4
$ws.Range("N74").Value2 = 2.4
$ws.Range("O74").Value2 = 4
$ws.Range("P74").Value2 = 2.6
$ws.Range("Q74").Value2 = 0
$ws.Range("R74").Value2 = 1.83
$ws.Range("S74").Value2 = 2.07
$ws.Range("T74").Value2 = 3.25
$ws.Range("U74").Value2 = 1.9
$ws.Range("V74").Value2 = 1.95
$ws.Range("W74").Value2 = 1.4
$ws.Range("X74").Value2 = -1
$ws.Range("Y74").Value2 = -1
$ws.Range("Z74").Value2 = 0.8300000000000001
$ws.Range("AA74").Value2 = -1
$ws.Range("AB74").Value2 = 0.8999999999999999
$ws.Range("AC74").Value2 = -1
$ws.Range("A104").Value2 = 102
$ws.Range("B104").Value2 = 7127374
$ws.Range("C104").Value2 = 'Australia ALeague'
$ws.Range("D104").Value2 = 'Australia ALeague'
$ws.Range("E104").Value2 = 45340.125
$ws.Range("F104").Value2 = 'Central Coast Mariners'
$ws.Range("G104").Value2 = 'Western Sydney Wanderers'
$ws.Range("H104").Value2 = 1
$ws.Range("I104").Value2 = 0
$ws.Range("J104").Value2 = 'H'
$ws.Range("K104").Value2 = 1.909
$ws.Range("L104").Value2 = 3.75
$ws.Range("M104").Value2 = 3.6
$ws.Range("N104").Value2 = 2.15
$ws.Range("O104").Value2 = 3.6
$ws.Range("P104").Value2 = 3.25
$ws.Range("Q104").Value2 = -0.25
$ws.Range("R104").Value2 = 1.86
$ws.Range("S104").Value2 = 2.04
$ws.Range("T104").Value2 = 2.75
$ws.Range("U104").Value2 = 1.975
$ws.Range("V104").Value2 = 1.875
$ws.Range("W104").Value2 = 1.15
$ws.Range("X104").Value2 = -1
$ws.Range("Y104").Value2 = -1
$ws.Range("Z104").Value2 = 0.8600000000000001
$ws.Range("AA104").Value2 = -1
$ws.Range("AB104").Value2 = -1
$ws.Range("AC104").Value2 = 0.875
$ws.Range("A105").Value2 = 103
$ws.Range("B105").Value2 = 7127370
$ws.Range("C105").Value2 = 'Australia ALeague'
$ws.Range("D105").Value2 = 'Australia ALeague'
$ws.Range("E105").Value2 = 45340.125
$ws.Range("F105").Value2 = 'Macarthur FC'
$ws.Range("G105").Value2 = 'Wellington Phoenix'
$ws.Range("H105").Value2 = 1
$ws.Range("I105").Value2 = 2
$ws.Range("J105").Value2 = 'A'
$ws.Range("K105").Value2 = 2.4
$ws.Range("L105").Value2 = 3.75
$ws.Range("M105").Value2 = 2.625
$ws.Range("N105").Value2 = 2.375
$ws.Range("O105").Value2 = 3.8
$ws.Range("P105").Value2 = 2.75
$ws.Range("Q105").Value2 = 0
$ws.Range("R105").Value2 = 1.8
$ws.Range("S105").Value2 = 2.05
$ws.Range("T105").Value2 = 3
$ws.Range("U105").Value2 = 1.9
$ws.Range("V105").Value2 = 1.95
$ws.Range("W105").Value2 = -1
$ws.Range("X105").Value2 = -1
$ws.Range("Y105").Value2 = 1.75
$ws.Range("Z105").Value2 = -1
$ws.Range("AA105").Value2 = 1.05
$ws.Range("AB105").Value2 = 0
$ws.Range("AC105").Value2 = -0
$ws.Range("A112").Value2 = 110
$ws.Range("B112").Value2 = 7127379
$ws.Range("C112").Value2 = 'Australia ALeague'
$ws.Range("D112").Value2 = 'Australia ALeague'
$ws.Range("E112").Value2 = 45347.125
$ws.Range("F112").Value2 = 'Melbourne Victory'
$ws.Range("G112").Value2 = 'Central Coast Mariners'
$ws.Range("H112").Value2 = 0
$ws.Range("I112").Value2 = 1
$ws.Range("J112").Value2 = 'A'
$ws.Range("K112").Value2 = 1.95
$ws.Range("L112").Value2 = 3.6
$ws.Range("M112").Value2 = 3.8
$ws.Range("N112").Value2 = 1.909
$ws.Range("O112").Value2 = 3.6
$ws.Range("P112").Value2 = 4
$ws.Range("Q112").Value2 = -0.5
$ws.Range("R112").Value2 = 1.9
$ws.Range("S112").Value2 = 1.95
$ws.Range("T112").Value2 = 2.75
$ws.Range("U112").Value2 = 1.925
$ws.Range("V112").Value2 = 1.925
$ws.Range("W112").Value2 = -1
$ws.Range("X112").Value2 = -1
$ws.Range("Y112").Value2 = 3
$ws.Range("Z112").Value2 = -1
$ws.Range("AA112").Value2 = 0.95
$ws.Range("AB112").Value2 = -1
$ws.Range("AC112").Value2 = 0.925
$ws.Range("A113").Value2 = 111
$ws.Range("B113").Value2 = 7127376
$ws.Range("C113").Value2 = 'Australia ALeague'
$ws.Range("D113").Value2 = 'Australia ALeague'
$ws.Range("E113").Value2 = 45347.125
$ws.Range("F113").Value2 = 'Newcastle Jets'
$ws.Range("G113").Value2 = 'Macarthur FC'
$ws.Range("H113").Value2 = 2
$ws.Range("I113").Value2 = 2
$ws.Range("J113").Value2 = 'D'
$ws.Range("K113").Value2 = 1.95
$ws.Range("L113").Value2 = 4
$ws.Range("M113").Value2 = 3.4
$ws.Range("N113").Value2 = 1.909
$ws.Range("O113").Value2 = 4.2
$ws.Range("P113").Value2 = 3.6
$ws.Range("Q113").Value2 = -0.5
$ws.Range("R113").Value2 = 1.89
$ws.Range("S113").Value2 = 2.01
$ws.Range("T113").Value2 = 3.5
$ws.Range("U113").Value2 = 1.95
$ws.Range("V113").Value2 = 1.9
$ws.Range("W113").Value2 = -1
$ws.Range("X113").Value2 = 3.2
$ws.Range("Y113").Value2 = -1
$ws.Range("Z113").Value2 = -1
$ws.Range("AA113").Value2 = 1.01
$ws.Range("AB113").Value2 = 0.95
$ws.Range("AC113").Value2 = -1
$ws.Range("A124").Value2 = 122
$ws.Range("B124").Value2 = 7128012
$ws.Range("C124").Value2 = 'Australia ALeague'
$ws.Range("D124").Value2 = 'Australia ALeague'
$ws.Range("E124").Value2 = 45361.125
$ws.Range("F124").Value2 = 'Macarthur FC'
$ws.Range("G124").Value2 = 'Central Coast Mariners'
$ws.Range("H124").Value2 = 0
$ws.Range("I124").Value2 = 3
$ws.Range("J124").Value2 = 'A'
$ws.Range("K124").Value2 = 2.4
$ws.Range("L124").Value2 = 3.5
$ws.Range("M124").Value2 = 2.75
$ws.Range("N124").Value2 = 3.4
$ws.Range("O124").Value2 = 3.75
$ws.Range("P124").Value2 = 2.05
$ws.Range("Q124").Value2 = 0.25
$ws.Range("R124").Value2 = 2.025
$ws.Range("S124").Value2 = 1.825
$ws.Range("T124").Value2 = 3
$ws.Range("U124").Value2 = 2.05
$ws.Range("V124").Value2 = 1.8
$ws.Range("W124").Value2 = -1
$ws.Range("X124").Value2 = -1
$ws.Range("Y124").Value2 = 1.05
$ws.Range("Z124").Value2 = -1
$ws.Range("AA124").Value2 = 0.825
$ws.Range("AB124").Value2 = 0
$ws.Range("AC124").Value2 = -0
$ws.Range("A125").Value2 = 123
$ws.Range("B125").Value2 = 7127388
$ws.Range("C125").Value2 = 'Australia ALeague'
$ws.Range("D125").Value2 = 'Australia ALeague'
$ws.Range("E125").Value2 = 45361.125
$ws.Range("F125").Value2 = 'Sydney FC'
$ws.Range("G125").Value2 = 'Brisbane Roar'
$ws.Range("H125").Value2 = 1
$ws.Range("I125").Value2 = 1
$ws.Range("J125").Value2 = 'D'
$ws.Range("K125").Value2 = 1.5
$ws.Range("L125").Value2 = 5
$ws.Range("M125").Value2 = 5
$ws.Range("N125").Value2 = 1.533
$ws.Range("O125").Value2 = 5.25
$ws.Range("P125").Value2 = 5
$ws.Range("Q125").Value2 = -1
$ws.Range("R125").Value2 = 1.8
$ws.Range("S125").Value2 = 2.05
$ws.Range("T125").Value2 = 3.5
$ws.Range("U125").Value2 = 1.925
$ws.Range("V125").Value2 = 1.925
$ws.Range("W125").Value2 = -1
$ws.Range("X125").Value2 = 4.25
$ws.Range("Y125").Value2 = -1
$ws.Range("Z125").Value2 = -1
$ws.Range("AA125").Value2 = 1.05
$ws.Range("AB125").Value2 = -1
$ws.Range("AC125").Value2 = 0.925
$ws.Range("A146").Value2 = 144
$ws.Range("B146").Value2 = 7127408
$ws.Range("C146").Value2 = 'Australia ALeague'
$ws.Range("D146").Value2 = 'Australia ALeague'
$ws.Range("E146").Value2 = 45395.10416666666
$ws.Range("F146").Value2 = 'Western United FC'
$ws.Range("G146").Value2 = 'Central Coast Mariners'
$ws.Range("K146").Value2 = 3.5
$ws.Range("L146").Value2 = 3.8
$ws.Range("M146").Value2 = 1.909
$ws.Range("N146").Value2 = 3.5
$ws.Range("O146").Value2 = 4.2
$ws.Range("P146").Value2 = 1.909
$ws.Range("Q146").Value2 = 0.5
$ws.Range("R146").Value2 = 1.99
$ws.Range("S146").Value2 = 1.91
$ws.Range("T146").Value2 = 3
$ws.Range("U146").Value2 = 1.85
$ws.Range("V146").Value2 = 2
$ws.Range("W146").Value2 = 0
$ws.Range("X146").Value2 = 0
$ws.Range("Y146").Value2 = 0
$ws.Range("Z146").Value2 = 0
$ws.Range("AA146").Value2 = 0
$ws.Range("A147").Value2 = 145
$ws.Range("B147").Value2 = 7127407
$ws.Range("C147").Value2 = 'Australia ALeague'
$ws.Range("D147").Value2 = 'Australia ALeague'
$ws.Range("E147").Value2 = 45395.1875
$ws.Range("F147").Value2 = 'Brisbane Roar'
$ws.Range("G147").Value2 = 'Newcastle Jets'
$ws.Range("K147").Value2 = 2.1
$ws.Range("L147").Value2 = 3.6
$ws.Range("M147").Value2 = 3.25
$ws.Range("N147").Value2 = 1.727
$ws.Range("O147").Value2 = 4.333
$ws.Range("P147").Value2 = 4.2
$ws.Range("Q147").Value2 = -0.75
$ws.Range("R147").Value2 = 1.89
$ws.Range("S147").Value2 = 2.01
$ws.Range("T147").Value2 = 3.25
$ws.Range("U147").Value2 = 1.825
$ws.Range("V147").Value2 = 2.025
$ws.Range("W147").Value2 = 0
$ws.Range("X147").Value2 = 0
$ws.Range("Y147").Value2 = 0
$ws.Range("Z147").Value2 = 0
$ws.Range("AA147").Value2 = 0
$ws.Range("A148").Value2 = 146
$ws.Range("B148").Value2 = 7127406
$ws.Range("C148").Value2 = 'Australia ALeague'
$ws.Range("D148").Value2 = 'Australia ALeague'
$ws.Range("E148").Value2 = 45395.28125
$ws.Range("F148").Value2 = 'Sydney FC'
$ws.Range("G148").Value2 = 'Western Sydney Wanderers'
$ws.Range("K148").Value2 = 1.833
$ws.Range("L148").Value2 = 4
$ws.Range("M148").Value2 = 3.75
$ws.Range("N148").Value2 = 1.65
$ws.Range("O148").Value2 = 4.5
$ws.Range("P148").Value2 = 4.5
$ws.Range("Q148").Value2 = -1
$ws.Range("R148").Value2 = 2.07
$ws.Range("S148").Value2 = 1.83
$ws.Range("T148").Value2 = 3.5
$ws.Range("U148").Value2 = 2
$ws.Range("V148").Value2 = 1.85
$ws.Range("W148").Value2 = 0
$ws.Range("X148").Value2 = 0
$ws.Range("Y148").Value2 = 0
$ws.Range("Z148").Value2 = 0
$ws.Range("AA148").Value2 = 0
$ws.Range("A149").Value2 = 147
$ws.Range("B149").Value2 = 7127409
$ws.Range("C149").Value2 = 'Australia ALeague'
$ws.Range("D149").Value2 = 'Australia ALeague'
$ws.Range("E149").Value2 = 45396.08333333334
$ws.Range("F149").Value2 = 'Melbourne City'
$ws.Range("G149").Value2 = 'Perth Glory'
$ws.Range("K149").Value2 = 1.571
$ws.Range("L149").Value2 = 4.5
$ws.Range("M149").Value2 = 4.75
$ws.Range("N149").Value2 = 1.4
$ws.Range("O149").Value2 = 5
$ws.Range("P149").Value2 = 7
$ws.Range("Q149").Value2 = -1.5
$ws.Range("R149").Value2 = 2.02
$ws.Range("S149").Value2 = 1.88
$ws.Range("T149").Value2 = 3.5
$ws.Range("U149").Value2 = 1.975
$ws.Range("V149").Value2 = 1.875
$ws.Range("W149").Value2 = 0
$ws.Range("X149").Value2 = 0
$ws.Range("Y149").Value2 = 0
$ws.Range("Z149").Value2 = 0
$ws.Range("AA149").Value2 = 0

$ws.Rows("150:151").Delete()
